$d = $word.ActiveDocument
Write-Host "Body content length:" $d.Content.End
$sec1 = $d.Sections(1)
$hdr1 = $sec1.Headers(1)
$hr = $hdr1.Range
Write-Host "Header StoryLength:" $hr.StoryLength
$hr.Find.Execute("LearningEDU") | Out-Null
Write-Host "Found:" $hr.Find.Found "Start:" $hr.Start "End:" $hr.End
